# Apply updated registration counts (Inscritos/Pagos/Inscrições homologadas)
# to the "Inscricoes" worksheet, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Map: row number -> hashtable of column letter -> new value
$changes = @{
    5  = @{ E = 2;  F = 2;          H = 2 }
    16 = @{ E = 9 }
    17 = @{ E = 91 }
    18 = @{ E = 94 }
    33 = @{ E = 30 }
    34 = @{         F = 6;          H = 6 }
    36 = @{ E = 81; F = 32;         H = 32 }
    41 = @{ E = 29 }
    46 = @{ E = 22 }
    48 = @{ E = 22; F = 12;         H = 12 }
    51 = @{ E = 6 }
    60 = @{ E = 15 }
    61 = @{ E = 24; F = 7;          H = 7 }
    63 = @{ E = 25; F = 7;          H = 7 }
    66 = @{ E = 29 }
    68 = @{ E = 13; F = 7;          H = 7 }
    72 = @{ E = 33; F = 16;         H = 16 }
    76 = @{ E = 41 }
    87 = @{ E = 12; F = 3;          H = 3 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
